$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.904.91'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.260.66'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.84'
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.18'
$ws.Range("E6").Value = '  +2.16%  '
$ws.Range("E7").Value = '  -1.04%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.489'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.99'
$ws.Range("E10").Value = '  +6.59%  '
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.64'
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.608.57'
$ws.Range("E14").Value = '  -0.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.38'
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.271.32'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.790'
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.796.13'
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("E19").Value = '  -3.44%  '
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.97'
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.02'
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.92'
$ws.Range("E23").Value = '  -3.02%  '
$ws.Range("E24").Value = '  -1.75%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E26").Value = '  -1.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.63'
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.49'
$ws.Range("E28").Value = '  +4.11%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.45'
$ws.Range("E29").Value = '  -2.48%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.10'
$ws.Range("E30").Value = '  +1.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '160.03'
$ws.Range("E31").Value = '  +0.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.22'
$ws.Range("E32").Value = '  -2.51%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.16'
$ws.Range("E34").Value = '  +3.77%  '
$ws.Range("E35").Value = '  -1.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '16.95'
$ws.Range("E36").Value = '  -1.75%  '
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.105'
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("E40").Value = '  -2.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.99'
$ws.Range("E41").Value = '  +0.66%  '
$ws.Range("E42").Value = '  +3.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.964.74'
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.57'
$ws.Range("E45").Value = '  -6.06%  '
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("E47").Value = '  -5.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.07'
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '72.71'
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("E50").Value = '  -1.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '90.98'
$ws.Range("E51").Value = '  -1.11%  '
